# Mark previously-TODO rubric items as done (Y) now that players can view
# their game records / bonus items have been implemented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: "客户端在等待匹配的过程中异常关闭..." -> mark as done
$ws.Cells.Item(17, 4).Value = "Y"

# Row 18: "当2个客户端都异常关闭，服务器应该解除该会话并记录log信息。" -> mark as done (new cell)
$ws.Cells.Item(18, 4).Value = "Y"

# Row 20: "账号管理系统..." bonus item -> mark as done
$ws.Cells.Item(20, 4).Value = "Y"

# Update the selected cell to reflect where the author was last working
$ws.Range("D20").Select()
